$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 19 new sponsor name/QQ-number pairs into columns L (name) / M (number),
# rows 47-65, extending the existing L/M list that previously stopped at row 45.
$ws.Cells.Item(47, 12).Value = "王老四"
$ws.Cells.Item(47, 13).Value = 904823142
$ws.Cells.Item(48, 12).Value = "欧豪666"
$ws.Cells.Item(48, 13).Value = 13683093
$ws.Cells.Item(49, 12).Value = "Biu凯凯"
$ws.Cells.Item(49, 13).Value = 1242297465
$ws.Cells.Item(50, 12).Value = "DDxZZD"
$ws.Cells.Item(50, 13).Value = 610812043
$ws.Cells.Item(51, 12).Value = "莫魂"
$ws.Cells.Item(51, 13).Value = 225975698
$ws.Cells.Item(52, 12).Value = "沈卡子"
$ws.Cells.Item(52, 13).Value = 1778330124
$ws.Cells.Item(53, 12).Value = "张思四"
$ws.Cells.Item(53, 13).Value = 1021570043
$ws.Cells.Item(54, 12).Value = "神天道至尊"
$ws.Cells.Item(54, 13).Value = 22278263
$ws.Cells.Item(55, 12).Value = "狂吃小笼包"
$ws.Cells.Item(55, 13).Value = 534339887
$ws.Cells.Item(56, 12).Value = "ghj999"
$ws.Cells.Item(56, 13).Value = 46380924
$ws.Cells.Item(57, 12).Value = "7kss"
$ws.Cells.Item(57, 13).Value = 434531878
$ws.Cells.Item(58, 12).Value = "阿迪王上"
$ws.Cells.Item(58, 13).Value = 48827206
$ws.Cells.Item(59, 12).Value = "灬Cyf"
$ws.Cells.Item(59, 13).Value = 23689807
$ws.Cells.Item(60, 12).Value = "奈奈吖丶"
$ws.Cells.Item(60, 13).Value = 601303810
$ws.Cells.Item(61, 12).Value = "丨朽木"
$ws.Cells.Item(61, 13).Value = 2024814813
$ws.Cells.Item(62, 12).Value = "小V亮1"
$ws.Cells.Item(62, 13).Value = 1965676132
$ws.Cells.Item(63, 12).Value = "逍-遥"
$ws.Cells.Item(63, 13).Value = 1356170852
$ws.Cells.Item(64, 12).Value = "小林先生0"
$ws.Cells.Item(64, 13).Value = 730864612
$ws.Cells.Item(65, 12).Value = "哔哔哔哔哔哔丶"
$ws.Cells.Item(65, 13).Value = 415945636

# Move the selection/viewport to match the post-edit state.
$ws.Range("O55").Select()
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 9
